$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 93 (shift existing rows 93:113 down to 94:114) ---
# A gap in the daily series (44234 -> 44236, missing 44235) is being backfilled.
$ws.Range("A93:D93").Insert(-4121)   # -4121 = xlShiftDown

# Copy the formatting from the row above (A92, styled as a date cell) onto the new A93
# so the new cell reuses the existing style instead of Excel creating a brand-new one.
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)  # -4122 = xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the values for the newly inserted row.
$ws.Range("A93").Value2 = 44235
$ws.Range("B93").Value2 = 0
$ws.Range("C93").Value2 = 0
$ws.Range("D93").Value2 = 0

# --- 2. The row that is now 112 (previously 111) gets its rolling-average columns filled in ---
$ws.Range("C112").Value2 = 1
$ws.Range("D112").Value2 = 43.78283712784589

# --- 3. Append a brand-new row (115) for the next day in the series (44257) ---
$ws.Range("A115").Value2 = 44257
$ws.Range("B115").Value2 = 0

# Match the styling Excel already uses for the date column by copying format from the row above.
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
